$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert the three new paragraphs ("2022.4.10:", the CDS_UPDATEREGISTRY
#    note, and a blank line) before the very first paragraph of the body.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$insertPoint = $firstPara.Range.Duplicate
$insertPoint.Collapse(1)

$newParasXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + `
          '<w:p>' + `
            '<w:pPr>' + `
              '<w:spacing w:line="220" w:lineRule="atLeast"/>' + `
              '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>' + `
            '</w:pPr>' + `
            '<w:r>' + `
              '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>' + `
              '<w:t>2022.4.10:</w:t>' + `
            '</w:r>' + `
          '</w:p>' + `
          '<w:p>' + `
            '<w:pPr>' + `
              '<w:spacing w:line="220" w:lineRule="atLeast"/>' + `
              '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>' + `
            '</w:pPr>' + `
            '<w:r>' + `
              '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>' + `
              '<w:t>' + [char]0x589E + [char]0x52A0 + '</w:t>' + `
            '</w:r>' + `
            '<w:r>' + `
              '<w:t>CDS_UPDATEREGISTRY</w:t>' + `
            '</w:r>' + `
            '<w:r>' + `
              '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>' + `
              '<w:t>' + [char]0x6807 + [char]0x5FD7 + [char]0xFF0C + [char]0x4F7F + [char]0x8BBE + [char]0x7F6E + [char]0x5168 + [char]0x5C40 + [char]0x751F + [char]0x6548 + [char]0x3002 + '</w:t>' + `
            '</w:r>' + `
          '</w:p>' + `
          '<w:p>' + `
            '<w:pPr>' + `
              '<w:spacing w:line="220" w:lineRule="atLeast"/>' + `
              '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>' + `
            '</w:pPr>' + `
          '</w:p>' + `
        '</w:body>' + `
      '</w:document>' + `
    '</pkg:xmlData>' + `
  '</pkg:part>' + `
'</pkg:package>'

$insertPoint.InsertXML($newParasXml)

# ---------------------------------------------------------------------------
# 2. The paragraph "代表顺时针旋转第一个扩展屏90度。" (now the 6th paragraph)
#    loses the <w:rPr> (rFonts hint=eastAsia) that used to sit on its <w:pPr>
#    paragraph mark, while keeping its three runs untouched.
# ---------------------------------------------------------------------------
$targetPara = $d.Paragraphs(6)

$replacementXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + `
          '<w:p>' + `
            '<w:pPr>' + `
              '<w:spacing w:line="220" w:lineRule="atLeast"/>' + `
            '</w:pPr>' + `
            '<w:r>' + `
              '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>' + `
              '<w:t>' + [char]0x4EE3 + [char]0x8868 + [char]0x987A + [char]0x65F6 + [char]0x9488 + [char]0x65CB + [char]0x8F6C + [char]0x7B2C + [char]0x4E00 + [char]0x4E2A + [char]0x6269 + [char]0x5C55 + [char]0x5C4F + '</w:t>' + `
            '</w:r>' + `
            '<w:r>' + `
              '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>' + `
              '<w:t>90</w:t>' + `
            '</w:r>' + `
            '<w:r>' + `
              '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>' + `
              '<w:t>' + [char]0x5EA6 + [char]0x3002 + '</w:t>' + `
            '</w:r>' + `
          '</w:p>' + `
        '</w:body>' + `
      '</w:document>' + `
    '</pkg:xmlData>' + `
  '</pkg:part>' + `
'</pkg:package>'

$targetPara.Range.InsertXML($replacementXml)
